# Relabel the per-category summary rows with "<Category>, <Metric>" labels,
# and collapse the "Summary" section (rows 34-40) into simplified
# "Total new nominations" / "Total carryover nominations" / "Total confirmed "
# / "Total unconfirmed " / "Total withdrawn " / "Total returned to the White
# House " rows (39 data rows total, down from 40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Civilian ---
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Confirmed "
$ws.Range("A9").Value  = "     Civilian, Unconfirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Returned to White House "

# --- Other Civilian ---
$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("A14").Value = "     Other Civilian, Confirmed "
$ws.Range("A15").Value = "     Other Civilian, Unconfirmed "
$ws.Range("A16").Value = "     Other Civilian, Withdrawn "

# --- Air Force ---
$ws.Range("A18").Value = "     Air Force, New nominations"
$ws.Range("A19").Value = "     Air Force, Confirmed "
$ws.Range("A20").Value = "     Air Force, Unconfirmed "

# --- Army ---
$ws.Range("A22").Value = "     Army, New nominations"
$ws.Range("A23").Value = "     Army, Confirmed "
$ws.Range("A24").Value = "     Army, Unconfirmed "
$ws.Range("A25").Value = "     Army, Returned to White House "

# --- Navy ---
$ws.Range("A27").Value = "     Navy, New nominations"
$ws.Range("A28").Value = "     Navy, Confirmed "
$ws.Range("A29").Value = "     Navy, Unconfirmed "

# --- Marine Corps ---
$ws.Range("A31").Value = "     Marine Corps, New nominations"
$ws.Range("A32").Value = "     Marine Corps, Confirmed "
$ws.Range("A33").Value = "     Marine Corps, Unconfirmed "

# --- Summary section: rewrite rows 34-39 in place, then drop old row 40 ---
$ws.Range("A34").Value = "Total new nominations"
$ws.Range("B34").Value = 24951
# Match the thousands-separator number format used by the other big totals
# (copy it in rather than setting NumberFormat directly so we reuse the
# existing style instead of registering a new one).
$ws.Range("B13").Copy()
$ws.Range("B34").PasteSpecial(-4122)

$ws.Range("A35").Value = "Total carryover nominations"

$ws.Range("A36").Value = "Total confirmed "
$ws.Range("B36").Value = 23050

$ws.Range("A37").Value = "Total unconfirmed "
$ws.Range("B37").Value = 1878

$ws.Range("A38").Value = "Total withdrawn "
$ws.Range("B38").Value = 15
# This cell used to hold a big "Total unconfirmed" number (thousands-sep
# style); the new value is small, so switch it back to the plain style
# used elsewhere for small totals.
$ws.Range("B7").Copy()
$ws.Range("B38").PasteSpecial(-4122)

$ws.Range("A39").Value = "Total returned to the White House "
$ws.Range("B39").Value = 8

# The old row 40 ("Total returned to the White House " / 8) is now redundant
# since its content was folded into row 39 above - remove it so the used
# range shrinks back down to A1:B39.
$ws.Rows("40").Delete()
